$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add final self-report BRIEF and CLASS responses in column D
$ws.Range("D2").Value = "Not stressful"
$ws.Range("D3").Value = "Very stressful"
$ws.Range("D4").Value = "Moderately stressful"
$ws.Range("D5").Value = "A little stressful"
$ws.Range("D6").Value = "Very stressful"
$ws.Range("D7").Value = "Very stressful"

# Widen column D to fit the new content
$ws.Columns.Item(4).ColumnWidth = 18

# Move the active cell selection down to D8 (below the new data)
[void]$ws.Range("D8").Select()
